$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 124 (shifts existing rows 124-236 down to 125-237).
$ws.Rows.Item(124).Insert()

# The newly-inserted row 124 is blank; populate it as a duplicate of the row that
# landed just below it (the former row 124, now at row 125), then overwrite the
# date with the new entry's date.
$ws.Range("A125:R125").Copy()
$ws.Range("A124").PasteSpecial()
$ws.Range("D124").Value = 44586
